# Add a new "yes" column (H is the spacer col already implied by new G values,
# the new data column itself is G/H) to the review table on Sheet1.
#
# Net effect (per the target diff):
#   - sharedStrings.xml gains one new shared string: "yes"
#   - Row 9 gets a new G9 = "yes" cell
#   - Rows 2-9 each gain a new, empty-but-styled H cell (same look as the
#     rest of the table: font "Mangal" 10pt, which is cellXfs style index 1)
#   - A brand new (otherwise empty) row 10 is added with styled G10/H10 cells
#   - Selection moves to G10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing populated rows are 2..9 (row 1 is the header). Stamp a new,
# empty H-column cell on every one of them so it picks up the same
# "Mangal" font styling already used throughout the sheet (this reuses the
# workbook's existing style slot instead of minting a new one).
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 8).Font.Name = "Mangal"
}

# Row 9: new G9 value "yes" plus the same styled-but-empty H9 cell.
$ws.Cells.Item(9, 7).Value = "yes"
$ws.Cells.Item(9, 7).Font.Name = "Mangal"
$ws.Cells.Item(9, 8).Font.Name = "Mangal"

# Brand new row 10: styled-but-empty G10 and H10 cells.
$ws.Cells.Item(10, 7).Font.Name = "Mangal"
$ws.Cells.Item(10, 8).Font.Name = "Mangal"

# Move the selection to match the post-edit cursor position.
$ws.Range("G10").Select()
